$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark after "Author"
$d.Bookmarks("_GoBack").Delete()

# 2. Fix "3-4 paragraphs" -> "3-4 sentences" typo
$d.Content.Find.Execute("of the reading in 3-4 paragraphs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "of the reading in 3-4 sentences", 2)

Write-Output "done"
